# first-train.xlsx: extend the data range from row 936 to row 1311 and
# update column C's recoded group value from 6 to 3 for the appended rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Recode column C for rows 937:1311 from 6 -> 3.
$ws.Range("C937:C1311").Value = 3

# 2) Update the _FilterDatabase defined name to cover the full data range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$C`$1:`$C`$1311"
    }
}

# 3) Re-apply the worksheet AutoFilter so its range covers C1:C1311
#    (toggle the existing one off first, since AutoFilter() on an
#    already-filtered range just flips it off).
$ws.Range("C1:C936").AutoFilter()
$ws.Range("C1:C1311").AutoFilter()

# 4) Scroll the view down so row 1036 is the top visible row, and select
#    C937:C1311 (the newly recoded block) as the active selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 1036
$win.ScrollColumn = 1
$ws.Range("C937:C1311").Select()
